$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 18:51"
$ws.Range("B4").Value = 6063589
$ws.Range("C4").Value = 16955
$ws.Range("D4").Value = 3352082
$ws.Range("E4").Value = 2526316
$ws.Range("G4").Value = 395
$ws.Range("H4").Value = 185191
$ws.Range("B5").Value = 3772945
$ws.Range("C5").Value = 8452
$ws.Range("E5").Value = 706707
$ws.Range("G5").Value = 262
$ws.Range("H5").Value = 118988
$ws.Range("B12").Value = 455621
$ws.Range("C12").Value = 3829
$ws.Range("G12").Value = 15
$ws.Range("H12").Value = 29011
$ws.Range("B20").Value = 265515
$ws.Range("C20").Value = 1517
$ws.Range("D20").Value = 241809
$ws.Range("E20").Value = 17461
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = 6245
$ws.Range("B21").Value = 265409
$ws.Range("C21").Value = 1462
$ws.Range("D21").Value = 206902
$ws.Range("E21").Value = 23035
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 35472
$ws.Range("B23").Value = 241609
$ws.Range("C23").Value = 1044
$ws.Range("E23").Value = 16754
$ws.Range("B24").Value = 223612
$ws.Range("C24").Value = 4177
$ws.Range("D24").Value = 164874
$ws.Range("E24").Value = 51924
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 6814
$ws.Range("B27").Value = 126971
$ws.Range("C27").Value = 123
$ws.Range("D27").Value = 112918
$ws.Range("E27").Value = 4945
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 9108
$ws.Range("B31").Value = 112000
$ws.Range("C31").Value = 1597
$ws.Range("D31").Value = 91051
$ws.Range("E31").Value = 20055
$ws.Range("G31").Value = 10
$ws.Range("H31").Value = 894
$ws.Range("D52").Value = 55337
$ws.Range("E52").Value = 1302
$ws.Range("B59").Value = 43403
$ws.Range("C59").Value = 387
$ws.Range("D59").Value = 30436
$ws.Range("E59").Value = 11484
$ws.Range("G59").Value = 8
$ws.Range("H59").Value = 1483
$ws.Range("A66").Value = "Azerbaiyan"
$ws.Range("B66").Value = 35986
$ws.Range("C66").Value = 142
$ws.Range("D66").Value = 33461
$ws.Range("E66").Value = 1998
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 527
$ws.Range("A67").Value = "Moldavia"
$ws.Range("B67").Value = 35904
$ws.Range("C67").Value = 358
$ws.Range("D67").Value = 24156
$ws.Range("E67").Value = 10767
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 981
$ws.Range("B74").Value = 23582
$ws.Range("C74").Value = 413
$ws.Range("D74").Value = 17447
$ws.Range("E74").Value = 5716
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 419
$ws.Range("A82").Value = "Libano"
$ws.Range("B82").Value = 15613
$ws.Range("C82").Value = 676
$ws.Range("D82").Value = 4260
$ws.Range("E82").Value = 11205
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = 148
$ws.Range("A83").Value = "Paraguay"
$ws.Range("B83").Value = 15290
$ws.Range("D83").Value = 8348
$ws.Range("E83").Value = 6677
$ws.Range("H83").Value = 265
$ws.Range("B90").Value = 10565
$ws.Range("C90").Value = 23
$ws.Range("E90").Value = 953
$ws.Range("B96").Value = 9195
$ws.Range("C96").Value = 112
$ws.Range("D96").Value = 4923
$ws.Range("E96").Value = 4001
$ws.Range("G96").Value = 5
$ws.Range("H96").Value = 271
$ws.Range("B141").Value = 1943
$ws.Range("C141").Value = 10
$ws.Range("D141").Value = 1107
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 563
$ws.Range("A144").Value = "Jordania"
$ws.Range("B144").Value = 1869
$ws.Range("C144").Value = 68
$ws.Range("D144").Value = 1367
$ws.Range("E144").Value = 487
$ws.Range("H144").Value = 15
$ws.Range("A145").Value = "Aruba"
$ws.Range("B145").Value = 1848
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 608
$ws.Range("E145").Value = 1232
$ws.Range("H145").Value = 8
$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 1820
$ws.Range("C146").Value = 32
$ws.Range("D146").Value = 1186
$ws.Range("E146").Value = 624
$ws.Range("H146").Value = 10
$ws.Range("A151").Value = "Republica de Chipre"
$ws.Range("B151").Value = 1481
$ws.Range("C151").Value = 14
$ws.Range("D151").Value = 935
$ws.Range("E151").Value = 526
$ws.Range("H151").Value = 20
$ws.Range("A152").Value = "Trinidad yTobago"
$ws.Range("B152").Value = 1476
$ws.Range("D152").Value = 594
$ws.Range("E152").Value = 867
$ws.Range("H152").Value = 15
$ws.Range("A212").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("C212").Value = 1
$ws.Range("D212").Value = 7
$ws.Range("E212").Value = 7
$ws.Range("A213").Value = "Groenlandia"
$ws.Range("B213").Value = 14
$ws.Range("D213").Value = 14
$ws.Range("E213").Value = 0
